# Refresh cryptos list: updates Price (D) and Volume/1h (E) columns for
# each coin row, and swaps the Aave / WEMIXTOKEN rows (50-51) to reflect
# their new ranking order, per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.408.49"
$ws.Range("E2").Value = "  -0.47%  "

# Row 3
$ws.Range("D3").Value = "2.100.03"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "

# Row 6
$ws.Range("E6").Value = "  +0.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.91%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.77%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08880"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "

# Row 13
$ws.Range("D13").Value = "2.100.15"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.787"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.005"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001145"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "  -0.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.276"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "

# Row 23
$ws.Range("D23").Value = "30.465.64"
$ws.Range("E23").Value = "  -0.42%  "

# Row 24
$ws.Range("E24").Value = "  +0.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.361"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "

# Row 26
$ws.Range("D26").Value = "2.345.64"
$ws.Range("E26").Value = "  -0.62%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.506"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.203"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1067"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.646"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.355"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.50%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.942"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.79%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.794"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02571"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06837"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2308"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.48%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6856"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.246"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.14%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.313"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6339"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.650"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "

# Row 48
$ws.Range("E48").Value = "  -0.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000341"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.28%  "

# Row 50
$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.12%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "
